$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J2").Value = "Homework 5"
$ws.Range("K2").Value = "Midterm 1"
